# Actualizacion a 4 de Abril.
# Adds daily COVID-19 case rows for 2020-04-01..03 (serials 43922-43924)
# to both "Hoja1" (regional case counts) and "Hoja2" (regional death counts).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# --- Hoja1: new rows 31-33 (columns A..S) ---
$hoja1Rows = @(
    @(43922, 30, 7, 12, 47, 3, 34, 156, 1636, 32, 82, 341, 275, 389, 68, 203, 5, 114, 3404),
    @(43923, 31, 12, 13, 54, 5, 35, 176, 1742, 33, 89, 370, 302, 432, 78, 234, 6, 156, 3737),
    @(43924, 32, 12, 17, 55, 5, 44, 185, 1957, 38, 97, 401, 347, 477, 89, 259, 6, 172, 4161)
)

$startRow = 31
for ($i = 0; $i -lt $hoja1Rows.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $hoja1Rows[$i]

    $dateCell = $ws1.Cells.Item($r, 1)
    $dateCell.Value = $rowVals[0]
    $dateCell.NumberFormat = "DD/MM/YY"

    for ($c = 1; $c -lt $rowVals.Length; $c++) {
        $ws1.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

# --- Hoja2: new rows 31-33 (columns A..S) ---
$hoja2Rows = @(
    @(43922, 30, 0, 0, 0, 0, 0, 1, 6, 0, 1, 0, 2, 7, 1, 0, 0, 0, 18),
    @(43923, 31, 0, 0, 0, 0, 0, 1, 8, 0, 1, 1, 2, 7, 1, 0, 0, 1, 22),
    @(43924, 32, 0, 0, 0, 0, 0, 1, 9, 0, 1, 2, 2, 9, 1, 0, 0, 2, 27)
)

for ($i = 0; $i -lt $hoja2Rows.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $hoja2Rows[$i]

    $dateCell = $ws2.Cells.Item($r, 1)
    $dateCell.Value = $rowVals[0]
    $dateCell.NumberFormat = "DD/MM/YY"

    for ($c = 1; $c -lt $rowVals.Length; $c++) {
        $ws2.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

# --- Selection / view state ---
# Hoja2: active cell moves to A31, no longer a multi-area I29:I30 + B1 selection.
$ws2.Activate()
$ws2.Range("A31").Select() | Out-Null

# Hoja1: stays the tab-selected sheet; active cell moves to B35 (view scrolled
# down near the new rows), replacing the old I29:I30 selection.
$ws1.Activate()
$ws1.Range("B35").Select() | Out-Null
